# TrialsSetup update (2026-01-22 12:00): the underlying Power Query refresh
# pulled updated "Days remaining" counts for three trials. Write the new
# values straight into the "Days remaining" column (column B) of Sheet1,
# matching what a live query refresh would have landed in the cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9  -> ALLEGRETTO-LTE (B7981028): 2  -> 1
$ws.Cells.Item(9, 2).Value = 1

# Row 11 -> REJOICE (MK-5909-003): 29 -> 28
$ws.Cells.Item(11, 2).Value = 28

# Row 14 -> REMASTER (CLOU): 49 -> 48
$ws.Cells.Item(14, 2).Value = 48
